$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert the new "trailing_ratio" column before the old "profit"
#    column (currently column Q). Everything from Q..V shifts right
#    by one (to R..W).
# ------------------------------------------------------------------
$ws.Columns("Q").Insert()
$ws.Range("Q1").Value = "trailing_ratio"
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 17).Value = 0.15   # column Q = 17
}

# ------------------------------------------------------------------
# 2) Insert the new "trailing_stop_limit_order_id" column before the
#    (now shifted) "timezone" column, which currently sits at V.
# ------------------------------------------------------------------
$ws.Columns("W").Insert()
$ws.Range("W1").Value = "trailing_stop_limit_order_id"

$trailingStopLimitOrderId = @{
    2  = "FAXXXX"
    3  = "FAXXXX"
    4  = "FAXXXX"
    5  = "FAXXXX"
    6  = "FAXXXX"
    7  = "FAXXXX"
    8  = "FAXXXX"
    9  = "FAXXXX"
    10 = "FAXXXX"
    11 = "FAXXXX"
    12 = "FAXXXX"
    13 = "FAXXXX"
}
foreach ($r in $trailingStopLimitOrderId.Keys) {
    $ws.Cells.Item($r, 23).Value = $trailingStopLimitOrderId[$r]   # column W = 23
}
# row 14 stays blank (no entry in the source data), but still touch the
# cell so a (empty) cell record is materialised for it, matching the
# other newly-inserted columns.
$ws.Cells.Item(14, 23).Style = "Normal"

# ------------------------------------------------------------------
# 3) Fill in the newly-available trailing_LIT_order_id values
#    (column V) for the rows that previously had it empty.
# ------------------------------------------------------------------
$trailingLitOrderId = @{
    10 = "FA195D1252483B2000"
    11 = "FA195D120DDD84A000"
    12 = "FA195D12895904A000"
    14 = "FA195D120A9BBB2000"
}
foreach ($r in $trailingLitOrderId.Keys) {
    $ws.Cells.Item($r, 22).Value = $trailingLitOrderId[$r]   # column V = 22
}

# ------------------------------------------------------------------
# 4) Renumber the trade "id" column (B) for rows 9-14.
# ------------------------------------------------------------------
$idFixups = @{
    9  = 7
    10 = 8
    11 = 9
    12 = 10
    13 = 11
    14 = 12
}
foreach ($r in $idFixups.Keys) {
    $ws.Cells.Item($r, 2).Value = $idFixups[$r]   # column B = 2
}
